$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title moved from 2023-01-19 to 2023-01-20
$ws.Range("A1").Value = "An Individual Customer's usage on 2023-01-20"

# Updated demand (kWh) readings in column B
$ws.Range("B3").Value  = 0
$ws.Range("B4").Value  = 0
$ws.Range("B5").Value  = 0
$ws.Range("B9").Value  = 0.0059983338
$ws.Range("B10").Value = 0.0000016662
$ws.Range("B13").Value = 0.005
$ws.Range("B14").Value = 0
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 0
$ws.Range("B20").Value = 0.0219958328
$ws.Range("B21").Value = 0.0200061111
$ws.Range("B22").Value = 0.005

# Last two hourly rows (20:00:00 and 21:00:00) no longer exist in the data
$ws.Range("A23:C24").EntireRow.Delete()
